$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the last data row (row 101) to make room for period 2508.
# This pushes the old last-data row (101) down to 102, and the signature block
# (rows 106-107) down to rows 107-108 -- matching the target layout.
$ws.Rows("101:101").Insert()

# Copy the formatting (borders/fonts/fills) of the row above into the newly
# inserted row so it matches the rest of the data table visually.
$ws.Range("B100:J100").Copy()
$ws.Range("B101").PasteSpecial(-4122)

# Update header summary figures.
$ws.Range("E11").Value = 2673625
$ws.Range("F13").Value = 87

# Rewrite the Periodo Mora / Valor Mora / Salario Basico table in chronological
# (ascending) order, 1705 .. 2508, one row per period.
$ws.Range("E16").Value = "1705"
$ws.Range("F16").Value = 9600
$ws.Range("G16").Value = 781242
$ws.Range("E17").Value = "1706"
$ws.Range("F17").Value = 9600
$ws.Range("G17").Value = 781242
$ws.Range("E18").Value = "1808"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 781242
$ws.Range("E19").Value = "1809"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242
$ws.Range("E20").Value = "1810"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242
$ws.Range("E21").Value = "1811"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 781242
$ws.Range("E22").Value = "1812"
$ws.Range("F22").Value = 31249
$ws.Range("G22").Value = 781242
$ws.Range("E23").Value = "1901"
$ws.Range("F23").Value = 31249
$ws.Range("G23").Value = 781242
$ws.Range("E24").Value = "1902"
$ws.Range("F24").Value = 31249
$ws.Range("G24").Value = 781242
$ws.Range("E25").Value = "1903"
$ws.Range("F25").Value = 31249
$ws.Range("G25").Value = 781242
$ws.Range("E26").Value = "1904"
$ws.Range("F26").Value = 31249
$ws.Range("G26").Value = 781242
$ws.Range("E27").Value = "1905"
$ws.Range("F27").Value = 31249
$ws.Range("G27").Value = 781242
$ws.Range("E28").Value = "1906"
$ws.Range("F28").Value = 31249
$ws.Range("G28").Value = 781242
$ws.Range("E29").Value = "1907"
$ws.Range("F29").Value = 31249
$ws.Range("G29").Value = 781242
$ws.Range("E30").Value = "1908"
$ws.Range("F30").Value = 31249
$ws.Range("G30").Value = 781242
$ws.Range("E31").Value = "1909"
$ws.Range("F31").Value = 31249
$ws.Range("G31").Value = 781242
$ws.Range("E32").Value = "1910"
$ws.Range("F32").Value = 31249
$ws.Range("G32").Value = 781242
$ws.Range("E33").Value = "1911"
$ws.Range("F33").Value = 31249
$ws.Range("G33").Value = 781242
$ws.Range("E34").Value = "1912"
$ws.Range("F34").Value = 31249
$ws.Range("G34").Value = 781242
$ws.Range("E35").Value = "2001"
$ws.Range("F35").Value = 31249
$ws.Range("G35").Value = 781242
$ws.Range("E36").Value = "2002"
$ws.Range("F36").Value = 31249
$ws.Range("G36").Value = 781242
$ws.Range("E37").Value = "2003"
$ws.Range("F37").Value = 31249
$ws.Range("G37").Value = 781242
$ws.Range("E38").Value = "2004"
$ws.Range("F38").Value = 31249
$ws.Range("G38").Value = 781242
$ws.Range("E39").Value = "2005"
$ws.Range("F39").Value = 31249
$ws.Range("G39").Value = 781242
$ws.Range("E40").Value = "2006"
$ws.Range("F40").Value = 31249
$ws.Range("G40").Value = 781242
$ws.Range("E41").Value = "2007"
$ws.Range("F41").Value = 31249
$ws.Range("G41").Value = 781242
$ws.Range("E42").Value = "2008"
$ws.Range("F42").Value = 31249
$ws.Range("G42").Value = 781242
$ws.Range("E43").Value = "2009"
$ws.Range("F43").Value = 31249
$ws.Range("G43").Value = 781242
$ws.Range("E44").Value = "2010"
$ws.Range("F44").Value = 31249
$ws.Range("G44").Value = 781242
$ws.Range("E45").Value = "2011"
$ws.Range("F45").Value = 31249
$ws.Range("G45").Value = 781242
$ws.Range("E46").Value = "2012"
$ws.Range("F46").Value = 31249
$ws.Range("G46").Value = 781242
$ws.Range("E47").Value = "2101"
$ws.Range("F47").Value = 31249
$ws.Range("G47").Value = 781242
$ws.Range("E48").Value = "2102"
$ws.Range("F48").Value = 31249
$ws.Range("G48").Value = 781242
$ws.Range("E49").Value = "2103"
$ws.Range("F49").Value = 31249
$ws.Range("G49").Value = 781242
$ws.Range("E50").Value = "2104"
$ws.Range("F50").Value = 31249
$ws.Range("G50").Value = 781242
$ws.Range("E51").Value = "2105"
$ws.Range("F51").Value = 31249
$ws.Range("G51").Value = 781242
$ws.Range("E52").Value = "2106"
$ws.Range("F52").Value = 31249
$ws.Range("G52").Value = 781242
$ws.Range("E53").Value = "2107"
$ws.Range("F53").Value = 31249
$ws.Range("G53").Value = 781242
$ws.Range("E54").Value = "2108"
$ws.Range("F54").Value = 31249
$ws.Range("G54").Value = 781242
$ws.Range("E55").Value = "2109"
$ws.Range("F55").Value = 31249
$ws.Range("G55").Value = 781242
$ws.Range("E56").Value = "2110"
$ws.Range("F56").Value = 31249
$ws.Range("G56").Value = 781242
$ws.Range("E57").Value = "2111"
$ws.Range("F57").Value = 31249
$ws.Range("G57").Value = 781242
$ws.Range("E58").Value = "2112"
$ws.Range("F58").Value = 31249
$ws.Range("G58").Value = 781242
$ws.Range("E59").Value = "2201"
$ws.Range("F59").Value = 31249
$ws.Range("G59").Value = 781242
$ws.Range("E60").Value = "2202"
$ws.Range("F60").Value = 31249
$ws.Range("G60").Value = 781242
$ws.Range("E61").Value = "2203"
$ws.Range("F61").Value = 31249
$ws.Range("G61").Value = 781242
$ws.Range("E62").Value = "2204"
$ws.Range("F62").Value = 31249
$ws.Range("G62").Value = 781242
$ws.Range("E63").Value = "2205"
$ws.Range("F63").Value = 31249
$ws.Range("G63").Value = 781242
$ws.Range("E64").Value = "2206"
$ws.Range("F64").Value = 31249
$ws.Range("G64").Value = 781242
$ws.Range("E65").Value = "2207"
$ws.Range("F65").Value = 31249
$ws.Range("G65").Value = 781242
$ws.Range("E66").Value = "2208"
$ws.Range("F66").Value = 31249
$ws.Range("G66").Value = 781242
$ws.Range("E67").Value = "2209"
$ws.Range("F67").Value = 31249
$ws.Range("G67").Value = 781242
$ws.Range("E68").Value = "2210"
$ws.Range("F68").Value = 31249
$ws.Range("G68").Value = 781242
$ws.Range("E69").Value = "2211"
$ws.Range("F69").Value = 31249
$ws.Range("G69").Value = 781242
$ws.Range("E70").Value = "2212"
$ws.Range("F70").Value = 31249
$ws.Range("G70").Value = 781242
$ws.Range("E71").Value = "2301"
$ws.Range("F71").Value = 31249
$ws.Range("G71").Value = 781242
$ws.Range("E72").Value = "2302"
$ws.Range("F72").Value = 31249
$ws.Range("G72").Value = 781242
$ws.Range("E73").Value = "2303"
$ws.Range("F73").Value = 31249
$ws.Range("G73").Value = 781242
$ws.Range("E74").Value = "2304"
$ws.Range("F74").Value = 31249
$ws.Range("G74").Value = 781242
$ws.Range("E75").Value = "2305"
$ws.Range("F75").Value = 31249
$ws.Range("G75").Value = 781242
$ws.Range("E76").Value = "2306"
$ws.Range("F76").Value = 31249
$ws.Range("G76").Value = 781242
$ws.Range("E77").Value = "2307"
$ws.Range("F77").Value = 31249
$ws.Range("G77").Value = 781242
$ws.Range("E78").Value = "2308"
$ws.Range("F78").Value = 31249
$ws.Range("G78").Value = 781242
$ws.Range("E79").Value = "2309"
$ws.Range("F79").Value = 31249
$ws.Range("G79").Value = 781242
$ws.Range("E80").Value = "2310"
$ws.Range("F80").Value = 31249
$ws.Range("G80").Value = 781242
$ws.Range("E81").Value = "2311"
$ws.Range("F81").Value = 31249
$ws.Range("G81").Value = 781242
$ws.Range("E82").Value = "2312"
$ws.Range("F82").Value = 31249
$ws.Range("G82").Value = 781242
$ws.Range("E83").Value = "2401"
$ws.Range("F83").Value = 31249
$ws.Range("G83").Value = 781242
$ws.Range("E84").Value = "2402"
$ws.Range("F84").Value = 31249
$ws.Range("G84").Value = 781242
$ws.Range("E85").Value = "2403"
$ws.Range("F85").Value = 31249
$ws.Range("G85").Value = 781242
$ws.Range("E86").Value = "2404"
$ws.Range("F86").Value = 31249
$ws.Range("G86").Value = 781242
$ws.Range("E87").Value = "2405"
$ws.Range("F87").Value = 31249
$ws.Range("G87").Value = 781242
$ws.Range("E88").Value = "2406"
$ws.Range("F88").Value = 31249
$ws.Range("G88").Value = 781242
$ws.Range("E89").Value = "2407"
$ws.Range("F89").Value = 31249
$ws.Range("G89").Value = 781242
$ws.Range("E90").Value = "2408"
$ws.Range("F90").Value = 31249
$ws.Range("G90").Value = 781242
$ws.Range("E91").Value = "2409"
$ws.Range("F91").Value = 31249
$ws.Range("G91").Value = 781242
$ws.Range("E92").Value = "2410"
$ws.Range("F92").Value = 31249
$ws.Range("G92").Value = 781242
$ws.Range("E93").Value = "2411"
$ws.Range("F93").Value = 31249
$ws.Range("G93").Value = 781242
$ws.Range("E94").Value = "2412"
$ws.Range("F94").Value = 31249
$ws.Range("G94").Value = 781242
$ws.Range("E95").Value = "2501"
$ws.Range("F95").Value = 31249
$ws.Range("G95").Value = 781242
$ws.Range("E96").Value = "2502"
$ws.Range("F96").Value = 31249
$ws.Range("G96").Value = 781242
$ws.Range("E97").Value = "2503"
$ws.Range("F97").Value = 31249
$ws.Range("G97").Value = 781242
$ws.Range("E98").Value = "2504"
$ws.Range("F98").Value = 31249
$ws.Range("G98").Value = 781242
$ws.Range("E99").Value = "2505"
$ws.Range("F99").Value = 31249
$ws.Range("G99").Value = 781242
$ws.Range("E100").Value = "2506"
$ws.Range("F100").Value = 31249
$ws.Range("G100").Value = 781242
$ws.Range("E101").Value = "2507"
$ws.Range("F101").Value = 31249
$ws.Range("G101").Value = 781242
$ws.Range("E102").Value = "2508"
$ws.Range("F102").Value = 31249
$ws.Range("G102").Value = 781242

# Shift the signature block text down to its new rows.
$ws.Range("B107").Value = "___________________________________"
$ws.Range("H107").Value = "___________________________________"
$ws.Range("B108").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H108").Value = "FIRMA DEL REPRESENTANTE LEGAL"

Write-Host "done"
